$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.545.74"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.755.57"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07471"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.087"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.76%  "

$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.004"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.149"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").Value = "1.752.77"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "

$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06402"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.736"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "27.587.19"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.067"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.15"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "1.953.42"
$ws.Range("E28").Value = "  +0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.082"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09217"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.661"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.529"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.24%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02285"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2093"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06016"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6281"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.923"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.185"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.784"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.718"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.939"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06901"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("E50").Value = "  -2.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "
